# "Corrección pts puche clasf"
# The classification table had wrong stats for "Puche": the commit corrects
# his Points (and the underlying G/P counts that produced them), which in
# turn changes his position once the table is re-sorted by Points
# (the sheet keeps an autoFilter + a "sort by Points desc" memory).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Correct "Puche"'s row (row 6: Puntos=B, Nºs1º=C, Max.pts=D, G=E, P=F, E=G, PJ=H) ---
# Puntos 9 -> 12
$ws.Range("B6").Value = 12
# G (goals/"ganados" column) 3 -> 4
$ws.Range("E6").Value = 4
# P ("perdidos" column) 1 -> 3
$ws.Range("F6").Value = 3

# --- 2. Re-sort the classification range A1:H7 (header in row 1) descending by Puntos (col B) ---
$tableRange = $ws.Range("A1:H7")
$tableRange.Sort($ws.Range("B1:B7"), 2, $null, $null, 1, $null, 1, 1)

# --- 3. (Re)apply the AutoFilter over the whole table ---
$tableRange.AutoFilter()

# Excel stores the autofilter range as a hidden sheet-scoped defined name;
# recreate it so the workbook round-trips the same way.
try {
    $filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$H`$7")
    $filterName.Visible = $false
} catch {
}

# --- 4. Update the active selection to cover the whole table, as left by the edit ---
$tableRange.Select()

Write-Host "Puche classification correction applied (9 -> 12 pts), table re-sorted and filtered."
